$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.757.93'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.75%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.224.25'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.86%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '298.84'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '84.25'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.514'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.07%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("E9").Value = '  -3.89%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0782'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '29.71'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.30'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -11.58%  '

$ws.Range("E13").Value = '  -2.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.567.47'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.77%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.29'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.27%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.09'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.60%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.233.37'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.91%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.718'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '39.669.52'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0877'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.05%  '

$ws.Range("E21").Value = '  -5.90%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.10'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.28%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.39'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.80%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '234.99'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.26%  '

$ws.Range("E25").Value = '  -0.07%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.42'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.17%  '

$ws.Range("E27").Value = '  -0.53%  '

$ws.Range("E28").Value = '  -2.85%  '

$ws.Range("E29").Value = '  -0.73%  '

$ws.Range("E30").Value = '  -1.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.25'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '149.74'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.71%  '

$ws.Range("E33").Value = '  -0.11%  '

$ws.Range("E34").Value = '  -5.65%  '

$ws.Range("E35").Value = '  -1.57%  '

$ws.Range("E36").Value = '  -2.48%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.38'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.38%  '

$ws.Range("E38").Value = '  -2.88%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0974'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.58%  '

$ws.Range("E40").Value = '  -4.17%  '

$ws.Range("E41").Value = '  -4.32%  '

$ws.Range("E42").Value = '  -5.49%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.926.14'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.31%  '

$ws.Range("E44").Value = '  -2.74%  '

$ws.Range("E45").Value = '  +0.79%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.23'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.74%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.50'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.47%  '

$ws.Range("E48").Value = '  -4.17%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.440.15'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.57%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '71.02'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '88.61'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.64%  '

